$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new BOM row (row 14) following the existing DNP/N-A pattern used by
# the other "not populated" rows (J1, J2, ISP1, BT1): a silkscreen/SAO
# header trace that isn't stuffed on this board revision.
$ws.Range("A14").Value = "SAO1"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = "SAO"
$ws.Range("D14").Value = "DNP"
$ws.Range("E14").Value = "N/A"

# Move the selection down to the next empty row, matching where the user's
# cursor ended up after typing the new row.
$ws.Range("A15").Select()
